# Helper: write $value into $cell while forcing text storage even when the
# string looks numeric (e.g. fund codes like "217024" or ratios like
# "35.05"). Plain `.Value =` would silently coerce those into numbers.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q1" sheet right before the "总计" sheet,
#    by duplicating an existing per-fund sheet (keeps the same sheetPr /
#    pageMargins / styles as its siblings) and overwriting its data.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Re-fetch, since inserting/copying sheets shifts positional references.
$totalSheet = $wb.Worksheets.Item("总计")

# The template sheet has 7 fund rows (rows 2-8); we only need 2, so drop
# the extra rows (this also shrinks the sheet dimension automatically).
$newSheet.Range("A4:H8").Clear()

# Row 2 data.
Set-TextValue $newSheet.Cells.Item(2, 2) "217024"
Set-TextValue $newSheet.Cells.Item(2, 3) "招商安盈债券"
Set-TextValue $newSheet.Cells.Item(2, 4) "35.05"
Set-TextValue $newSheet.Cells.Item(2, 5) "20.20"
Set-TextValue $newSheet.Cells.Item(2, 6) "0.91"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.3190"
$newSheet.Cells.Item(2, 8).Value = 5

# Row 3 data.
Set-TextValue $newSheet.Cells.Item(3, 2) "350002"
Set-TextValue $newSheet.Cells.Item(3, 3) "天治低碳经济灵活配置混合"
Set-TextValue $newSheet.Cells.Item(3, 4) "0.76"
Set-TextValue $newSheet.Cells.Item(3, 5) "65.23"
Set-TextValue $newSheet.Cells.Item(3, 6) "6.61"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.0502"
$newSheet.Cells.Item(3, 8).Value = 2

# ------------------------------------------------------------------
# 2. Add the new "2022-Q1" summary row at the top of the "总计" sheet,
#    pushing every existing row down by one.
# ------------------------------------------------------------------

# Shift the 5 existing data rows (2..6) down to (3..7), bottom-up so
# each source row is read before it gets overwritten. `.Text` is used
# for the read because `.Value` does not resolve reliably here; writing
# the resulting numeric-looking strings back through `.Value` restores
# the proper numeric storage for columns C/D.
for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $bVal = $totalSheet.Cells.Item($r, 2).Text
    $cVal = $totalSheet.Cells.Item($r, 3).Text
    $dVal = $totalSheet.Cells.Item($r, 4).Text
    $totalSheet.Cells.Item($dst, 2).Value = $bVal
    $totalSheet.Cells.Item($dst, 3).Value = $cVal
    $totalSheet.Cells.Item($dst, 4).Value = $dVal
}

# Row 7 (2020-Q4, the former last row) needs a fresh A7 index cell --
# it did not exist before the shift. Clone the styling used by the rest
# of the index column (A2:A6) and set its value.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$totalSheet.Cells.Item(7, 1).Value = 5

# New first data row: 2022-Q1.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.37

# Re-sequence the index column for the rows that moved down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4

# Restore the originally-active sheet/tab (sheet-juggling above leaves the
# newest sheet selected, which the source workbook never had).
$wb.Worksheets.Item("2020-Q4").Activate()
